$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Reference Code" row (row 2) to "Location Number" to reflect the
# new active/inactive report state.
$ws.Range("A2:D2").Value = "Location Number"

# Move the active cell/selection to D2 as reflected in the saved view state.
$ws.Range("D2").Select()
